$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4941.6
$ws.Range("I15").Value = 4941.6
$ws.Range("K15").Value = 14824.8
$ws.Range("M15").Value = -14655.8
$ws.Range("H45").Value = 6388.143
$ws.Range("J45").Value = 7780
$ws.Range("L45").Value = 23340
$ws.Range("N45").Value = -23724
$ws.Range("H52").Value = 1027.7142
$ws.Range("I52").Value = 199
$ws.Range("J52").Value = 6000
$ws.Range("K52").Value = 597
$ws.Range("L52").Value = 18000
$ws.Range("M52").Value = -437
$ws.Range("N52").Value = -18320
$ws.Range("H80").Value = 13889820
$ws.Range("I80").Value = 20834032
$ws.Range("J80").Value = 1398.625
$ws.Range("K80").Value = 62502096
$ws.Range("L80").Value = 4195.875
$ws.Range("M80").Value = -62501098
$ws.Range("N80").Value = -6191.875
$ws.Range("H83").Value = 13889820
$ws.Range("I83").Value = 20834032
$ws.Range("J83").Value = 1398.625
$ws.Range("K83").Value = 187506288
$ws.Range("L83").Value = 12587.625
$ws.Range("M83").Value = -187501296
$ws.Range("N83").Value = -22571.625
$ws.Range("H86").Value = 2614.9
$ws.Range("I86").Value = 2631.1875
$ws.Range("J86").Value = 2549.75
$ws.Range("K86").Value = 2631.1875
$ws.Range("L86").Value = 2549.75
$ws.Range("M86").Value = -1508.1875
$ws.Range("N86").Value = -4795.75
$ws.Range("H89").Value = 2614.9
$ws.Range("I89").Value = 2631.1875
$ws.Range("J89").Value = 2549.75
$ws.Range("K89").Value = 13155.9375
$ws.Range("L89").Value = 12748.75
$ws.Range("M89").Value = -7539.9375
$ws.Range("N89").Value = -23980.75
$ws.Range("H109").Value = 40005.332
$ws.Range("J109").Value = 40005.332
$ws.Range("L109").Value = 40005.332
$ws.Range("N109").Value = -42779.332
$ws.Range("H113").Value = 2147.7173
$ws.Range("I113").Value = 1865
$ws.Range("J113").Value = 2198.4614
$ws.Range("K113").Value = 1865
$ws.Range("L113").Value = 2198.4614
$ws.Range("M113").Value = 1389
$ws.Range("N113").Value = -8706.4614
$ws.Range("H114").Value = 45941.6
$ws.Range("J114").Value = 45941.6
$ws.Range("L114").Value = 45941.6
$ws.Range("N114").Value = -54619.6
$ws.Range("H117").Value = 46997.332
$ws.Range("J117").Value = 46997.332
$ws.Range("L117").Value = 46997.332
$ws.Range("N117").Value = -56175.332
$ws.Range("H120").Value = 49150.668
$ws.Range("J120").Value = 49150.668
$ws.Range("L120").Value = 49150.668
$ws.Range("N120").Value = -58826.668
$ws.Range("H123").Value = 37400
$ws.Range("J123").Value = 37400
$ws.Range("L123").Value = 37400
$ws.Range("N123").Value = -47200
$ws.Range("H124").Value = 40137.75
$ws.Range("J124").Value = 48517
$ws.Range("L124").Value = 48517
$ws.Range("N124").Value = -58337
$ws.Range("H126").Value = 43753
$ws.Range("J126").Value = 43753
$ws.Range("L126").Value = 43753
$ws.Range("N126").Value = -53633
$ws.Range("H130").Value = 40328
$ws.Range("J130").Value = 40328
$ws.Range("L130").Value = 40328
$ws.Range("N130").Value = -50368
$ws.Range("H133").Value = 46754.145
$ws.Range("J133").Value = 46754.145
$ws.Range("L133").Value = 46754.145
$ws.Range("N133").Value = -56874.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7981.067
$ws.Range("I32").Value = 7091.525
$ws.Range("K32").Value = 7091.525
$ws.Range("M32").Value = -6804.525
$ws.Range("H37").Value = 49000
$ws.Range("I37").Value = 50000
$ws.Range("J37").Value = 48000
$ws.Range("K37").Value = 50000
$ws.Range("L37").Value = 48000
$ws.Range("M37").Value = -49727
$ws.Range("N37").Value = -48546
$ws.Range("H44").Value = 35234.215
$ws.Range("J44").Value = 35234.215
$ws.Range("L44").Value = 35234.215
$ws.Range("N44").Value = -36210.215
$ws.Range("H61").Value = 1465.26
$ws.Range("I61").Value = 1261.5111
$ws.Range("J61").Value = 3299
$ws.Range("K61").Value = 1261.5111
$ws.Range("L61").Value = 3299
$ws.Range("M61").Value = -1049.5111
$ws.Range("N61").Value = -3723
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H76").Value = 21920
$ws.Range("J76").Value = 21920
$ws.Range("L76").Value = 21920
$ws.Range("N76").Value = -22596
$ws.Range("H79").Value = 21920
$ws.Range("J79").Value = 21920
$ws.Range("L79").Value = 21920
$ws.Range("N79").Value = -24260
$ws.Range("H80").Value = 48534.125
$ws.Range("J80").Value = 48534.125
$ws.Range("L80").Value = 48534.125
$ws.Range("N80").Value = -50530.125
$ws.Range("H83").Value = 48534.125
$ws.Range("J83").Value = 48534.125
$ws.Range("L83").Value = 145602.375
$ws.Range("N83").Value = -155586.375
$ws.Range("H109").Value = 45666.2
$ws.Range("J109").Value = 45666.2
$ws.Range("L109").Value = 45666.2
$ws.Range("N109").Value = -48440.2
$ws.Range("H117").Value = 43573.4
$ws.Range("J117").Value = 43573.4
$ws.Range("L117").Value = 43573.4
$ws.Range("N117").Value = -52751.4
$ws.Range("H118").Value = 48387.332
$ws.Range("J118").Value = 48387.332
$ws.Range("L118").Value = 48387.332
$ws.Range("N118").Value = -51701.332
$ws.Range("H119").Value = 49796
$ws.Range("J119").Value = 49796
$ws.Range("L119").Value = 49796
$ws.Range("N119").Value = -59472
$ws.Range("H121").Value = 43471
$ws.Range("J121").Value = 43471
$ws.Range("L121").Value = 43471
$ws.Range("N121").Value = -46965
$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381
$ws.Range("H131").Value = 48041
$ws.Range("J131").Value = 48041
$ws.Range("L131").Value = 48041
$ws.Range("N131").Value = -58121
$ws.Range("H132").Value = 9617045
$ws.Range("I132").Value = 12501003
$ws.Range("J132").Value = 3851.25
$ws.Range("K132").Value = 37503009
$ws.Range("L132").Value = 11553.75
$ws.Range("M132").Value = -37500479
$ws.Range("N132").Value = -16613.75
$ws.Range("H133").Value = 40532.2
$ws.Range("J133").Value = 40532.2
$ws.Range("L133").Value = 40532.2
$ws.Range("N133").Value = -45592.2
$ws.Range("H134").Value = 45690
$ws.Range("J134").Value = 45690
$ws.Range("L134").Value = 45690
$ws.Range("N134").Value = -55830
$ws.Range("H136").Value = 1465.26
$ws.Range("I136").Value = 1261.5111
$ws.Range("J136").Value = 3299
$ws.Range("K136").Value = 3784.5333
$ws.Range("L136").Value = 9897
$ws.Range("M136").Value = -1234.5333
$ws.Range("N136").Value = -14997
$ws.Range("H137").Value = 42950
$ws.Range("J137").Value = 42950
$ws.Range("L137").Value = 42950
$ws.Range("N137").Value = -53150
$ws.Range("H138").Value = 47000
$ws.Range("J138").Value = 47000
$ws.Range("L138").Value = 47000
$ws.Range("N138").Value = -57280
$ws.Range("H139").Value = 45944.855
$ws.Range("J139").Value = 45944.855
$ws.Range("L139").Value = 45944.855
$ws.Range("N139").Value = -56224.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 54999.5
$ws.Range("J57").Value = 54999.5
$ws.Range("L57").Value = 54999.5
$ws.Range("N57").Value = -56439.5
$ws.Range("H94").Value = 900.4286
$ws.Range("I94").Value = 979.5
$ws.Range("J94").Value = 795
$ws.Range("K94").Value = 979.5
$ws.Range("L94").Value = 795
$ws.Range("M94").Value = -528.5
$ws.Range("N94").Value = -1697
$ws.Range("H107").Value = 2199.4138
$ws.Range("I107").Value = 2132.7368
$ws.Range("J107").Value = 2326.1
$ws.Range("K107").Value = 2132.7368
$ws.Range("L107").Value = 2326.1
$ws.Range("M107").Value = -212.7368000000001
$ws.Range("N107").Value = -6166.1
$ws.Range("H108").Value = 44184
$ws.Range("J108").Value = 44184
$ws.Range("L108").Value = 44184
$ws.Range("N108").Value = -51864
$ws.Range("H110").Value = 48702
$ws.Range("J110").Value = 48702
$ws.Range("L110").Value = 48702
$ws.Range("N110").Value = -56882
$ws.Range("H112").Value = 46665
$ws.Range("J112").Value = 46665
$ws.Range("L112").Value = 46665
$ws.Range("N112").Value = -49619
$ws.Range("H116").Value = 43538
$ws.Range("J116").Value = 43538
$ws.Range("L116").Value = 43538
$ws.Range("N116").Value = -52716
$ws.Range("H117").Value = 48734
$ws.Range("J117").Value = 48734
$ws.Range("L117").Value = 48734
$ws.Range("N117").Value = -57912
$ws.Range("H119").Value = 46876.5
$ws.Range("J119").Value = 46876.5
$ws.Range("L119").Value = 46876.5
$ws.Range("N119").Value = -56552.5
$ws.Range("H120").Value = 47761
$ws.Range("J120").Value = 47761
$ws.Range("L120").Value = 47761
$ws.Range("N120").Value = -57437
$ws.Range("H124").Value = 52496
$ws.Range("J124").Value = 52496
$ws.Range("L124").Value = 52496
$ws.Range("N124").Value = -62316
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652
$ws.Range("H130").Value = 49387.5
$ws.Range("J130").Value = 49387.5
$ws.Range("L130").Value = 49387.5
$ws.Range("N130").Value = -59427.5
$ws.Range("H132").Value = 44225
$ws.Range("J132").Value = 44225
$ws.Range("L132").Value = 44225
$ws.Range("N132").Value = -54345
$ws.Range("H133").Value = 44666.332
$ws.Range("J133").Value = 44666.332
$ws.Range("L133").Value = 44666.332
$ws.Range("N133").Value = -54786.332
$ws.Range("H134").Value = 3504.5898
$ws.Range("I134").Value = 1752.4762
$ws.Range("J134").Value = 4150.1055
$ws.Range("K134").Value = 5257.4286
$ws.Range("L134").Value = 12450.3165
$ws.Range("M134").Value = -2722.4286
$ws.Range("N134").Value = -17520.3165
$ws.Range("H136").Value = 54999.5
$ws.Range("J136").Value = 54999.5
$ws.Range("L136").Value = 54999.5
$ws.Range("N136").Value = -65199.5
$ws.Range("H139").Value = 43837.6
$ws.Range("I139").Value = 10709
$ws.Range("J139").Value = 52119.75
$ws.Range("K139").Value = 10709
$ws.Range("L139").Value = 52119.75
$ws.Range("M139").Value = -5569
$ws.Range("N139").Value = -62399.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225.18
$ws.Range("I31").Value = 1010.4889
$ws.Range("J31").Value = 3219.018
$ws.Range("K31").Value = 1010.4889
$ws.Range("L31").Value = 3219.018
$ws.Range("M31").Value = -715.4889
$ws.Range("N31").Value = -3809.018
$ws.Range("H34").Value = 2225.18
$ws.Range("I34").Value = 1010.4889
$ws.Range("J34").Value = 3219.018
$ws.Range("K34").Value = 1010.4889
$ws.Range("L34").Value = 3219.018
$ws.Range("M34").Value = -808.4889
$ws.Range("N34").Value = -3623.018
$ws.Range("H52").Value = 32997.5
$ws.Range("J52").Value = 32997.5
$ws.Range("L52").Value = 32997.5
$ws.Range("N52").Value = -33585.5
$ws.Range("H99").Value = 3339.818
$ws.Range("J99").Value = 3412.6667
$ws.Range("L99").Value = 3412.6667
$ws.Range("N99").Value = -6408.6667
$ws.Range("H110").Value = 38898.6
$ws.Range("J110").Value = 38898.6
$ws.Range("L110").Value = 38898.6
$ws.Range("N110").Value = -47078.6
$ws.Range("H112").Value = 42698
$ws.Range("J112").Value = 42698
$ws.Range("L112").Value = 42698
$ws.Range("N112").Value = -45652
$ws.Range("H116").Value = 49368.5
$ws.Range("J116").Value = 49368.5
$ws.Range("L116").Value = 49368.5
$ws.Range("N116").Value = -58546.5
$ws.Range("H119").Value = 46573.75
$ws.Range("J119").Value = 46573.75
$ws.Range("L119").Value = 46573.75
$ws.Range("N119").Value = -56249.75
$ws.Range("H122").Value = 201106.17
$ws.Range("I122").Value = 600506
$ws.Range("J122").Value = 1406.25
$ws.Range("K122").Value = 1801518
$ws.Range("L122").Value = 4218.75
$ws.Range("M122").Value = -1799068
$ws.Range("N122").Value = -9118.75
$ws.Range("H126").Value = 3339.818
$ws.Range("J126").Value = 3412.6667
$ws.Range("K126").Value = 9036
$ws.Range("L126").Value = 10238.0001
$ws.Range("N126").Value = -15178.0001
$ws.Range("H132").Value = 47327.547
$ws.Range("I132").Value = 1954.0741
$ws.Range("J132").Value = 353598.5
$ws.Range("K132").Value = 5862.2223
$ws.Range("L132").Value = 1060795.5
$ws.Range("M132").Value = -3332.2223
$ws.Range("N132").Value = -1065855.5
$ws.Range("H133").Value = 27849.6
$ws.Range("J133").Value = 27849.6
$ws.Range("L133").Value = 27849.6
$ws.Range("N133").Value = -32909.6
$ws.Range("H137").Value = 43474.75
$ws.Range("J137").Value = 43474.75
$ws.Range("L137").Value = 43474.75
$ws.Range("N137").Value = -53674.75
$ws.Range("H139").Value = 57855.8
$ws.Range("I139").Value = 42000
$ws.Range("J139").Value = 61819.75
$ws.Range("K139").Value = 42000
$ws.Range("L139").Value = 61819.75
$ws.Range("N139").Value = -72099.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 328.42856
$ws.Range("I26").Value = 349.75
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 1049.25
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = -761.25
$ws.Range("N26").Value = -1476
$ws.Range("H38").Value = 125335.375
$ws.Range("I38").Value = 300
$ws.Range("J38").Value = 167013.83
$ws.Range("K38").Value = 900
$ws.Range("L38").Value = 501041.49
$ws.Range("M38").Value = -553
$ws.Range("N38").Value = -501735.49
$ws.Range("H64").Value = 2826.7222
$ws.Range("I64").Value = 1522.5
$ws.Range("J64").Value = 3199.3572
$ws.Range("K64").Value = 4567.5
$ws.Range("L64").Value = 9598.0716
$ws.Range("M64").Value = -4297.5
$ws.Range("N64").Value = -10138.0716
$ws.Range("H67").Value = 2826.7222
$ws.Range("I67").Value = 1522.5
$ws.Range("J67").Value = 3199.3572
$ws.Range("K67").Value = 4567.5
$ws.Range("L67").Value = 9598.0716
$ws.Range("M67").Value = -3631.5
$ws.Range("N67").Value = -11470.0716
$ws.Range("H76").Value = 4537.5186
$ws.Range("I76").Value = 1506.5
$ws.Range("K76").Value = 4519.5
$ws.Range("M76").Value = -4136.5
$ws.Range("H79").Value = 4537.5186
$ws.Range("I79").Value = 1506.5
$ws.Range("K79").Value = 4519.5
$ws.Range("M79").Value = -3193.5
$ws.Range("H104").Value = 3142.8572
$ws.Range("J104").Value = 3142.8572
$ws.Range("L104").Value = 9428.5716
$ws.Range("N104").Value = -14670.5716
$ws.Range("H131").Value = 3677.3171
$ws.Range("I131").Value = 12873.375
$ws.Range("J131").Value = 1447.9697
$ws.Range("K131").Value = 38620.125
$ws.Range("L131").Value = 4343.909100000001
$ws.Range("M131").Value = -33580.125
$ws.Range("N131").Value = -14423.9091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47659.332
$ws.Range("J110").Value = 47659.332
$ws.Range("L110").Value = 47659.332
$ws.Range("N110").Value = -55839.332
$ws.Range("H116").Value = 37798.4
$ws.Range("J116").Value = 37798.4
$ws.Range("L116").Value = 37798.4
$ws.Range("N116").Value = -46976.4
$ws.Range("H122").Value = 1128.0714
$ws.Range("I122").Value = 1129.7
$ws.Range("J122").Value = 1124
$ws.Range("K122").Value = 3389.1
$ws.Range("L122").Value = 3372
$ws.Range("M122").Value = -939.1000000000004
$ws.Range("N122").Value = -8272
$ws.Range("H124").Value = 39500
$ws.Range("J124").Value = 39500
$ws.Range("L124").Value = 39500
$ws.Range("N124").Value = -49320
$ws.Range("H130").Value = 44891.332
$ws.Range("J130").Value = 44891.332
$ws.Range("L130").Value = 44891.332
$ws.Range("N130").Value = -54931.332
$ws.Range("H132").Value = 2836.276
$ws.Range("I132").Value = 1869.6
$ws.Range("J132").Value = 3872
$ws.Range("K132").Value = 5608.799999999999
$ws.Range("L132").Value = 11616
$ws.Range("M132").Value = -3078.799999999999
$ws.Range("N132").Value = -16676
$ws.Range("H135").Value = 39933.332
$ws.Range("J135").Value = 39933.332
$ws.Range("L135").Value = 39933.332
$ws.Range("N135").Value = -50073.332
$ws.Range("H137").Value = 39800
$ws.Range("J137").Value = 39800
$ws.Range("L137").Value = 39800
$ws.Range("N137").Value = -50000
$ws.Range("H138").Value = 42000
$ws.Range("J138").Value = 42000
$ws.Range("L138").Value = 42000
$ws.Range("N138").Value = -52280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2126.7742
$ws.Range("I7").Value = 1852.7778
$ws.Range("J7").Value = 3976.25
$ws.Range("K7").Value = 1852.7778
$ws.Range("L7").Value = 3976.25
$ws.Range("M7").Value = -1740.7778
$ws.Range("N7").Value = -4200.25
$ws.Range("H61").Value = 6777
$ws.Range("I61").Value = 7250
$ws.Range("J61").Value = 6461.6665
$ws.Range("K61").Value = 7250
$ws.Range("L61").Value = 6461.6665
$ws.Range("M61").Value = -7048
$ws.Range("N61").Value = -6865.6665
$ws.Range("H108").Value = 48618
$ws.Range("J108").Value = 48618
$ws.Range("L108").Value = 48618
$ws.Range("N108").Value = -56298
$ws.Range("H109").Value = 35273
$ws.Range("J109").Value = 35273
$ws.Range("L109").Value = 35273
$ws.Range("N109").Value = -38047
$ws.Range("H110").Value = 44570
$ws.Range("J110").Value = 44570
$ws.Range("L110").Value = 44570
$ws.Range("N110").Value = -52750
$ws.Range("H113").Value = 6777
$ws.Range("I113").Value = 7250
$ws.Range("J113").Value = 6461.6665
$ws.Range("K113").Value = 7250
$ws.Range("L113").Value = 6461.6665
$ws.Range("M113").Value = -5080
$ws.Range("N113").Value = -10801.6665
$ws.Range("H114").Value = 39394
$ws.Range("J114").Value = 39394
$ws.Range("L114").Value = 39394
$ws.Range("N114").Value = -48072
$ws.Range("H116").Value = 47992
$ws.Range("J116").Value = 47992
$ws.Range("L116").Value = 47992
$ws.Range("N116").Value = -57170
$ws.Range("H117").Value = 43380
$ws.Range("J117").Value = 43380
$ws.Range("L117").Value = 43380
$ws.Range("N117").Value = -52558
$ws.Range("H118").Value = 40350.25
$ws.Range("J118").Value = 40350.25
$ws.Range("L118").Value = 40350.25
$ws.Range("N118").Value = -43664.25
$ws.Range("H119").Value = 45097.332
$ws.Range("J119").Value = 45097.332
$ws.Range("L119").Value = 45097.332
$ws.Range("N119").Value = -54773.332
$ws.Range("H120").Value = 50991.5
$ws.Range("J120").Value = 50991.5
$ws.Range("L120").Value = 50991.5
$ws.Range("N120").Value = -60667.5
$ws.Range("H122").Value = 202181.8
$ws.Range("I122").Value = 251976
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 755928
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -753478
$ws.Range("N122").Value = -13915
$ws.Range("H125").Value = 45401.75
$ws.Range("J125").Value = 45401.75
$ws.Range("L125").Value = 45401.75
$ws.Range("N125").Value = -55241.75
$ws.Range("H126").Value = 2126.7742
$ws.Range("I126").Value = 1852.7778
$ws.Range("J126").Value = 3976.25
$ws.Range("K126").Value = 5558.3334
$ws.Range("L126").Value = 11928.75
$ws.Range("M126").Value = -3088.3334
$ws.Range("N126").Value = -16868.75
$ws.Range("H128").Value = 47992
$ws.Range("J128").Value = 47992
$ws.Range("L128").Value = 47992
$ws.Range("N128").Value = -57952
$ws.Range("H133").Value = 34998.4
$ws.Range("J133").Value = 34998.4
$ws.Range("L133").Value = 34998.4
$ws.Range("N133").Value = -40058.4
$ws.Range("H134").Value = 44054.832
$ws.Range("J134").Value = 44054.832
$ws.Range("L134").Value = 44054.832
$ws.Range("N134").Value = -54194.832
$ws.Range("H137").Value = 34975
$ws.Range("J137").Value = 34975
$ws.Range("L137").Value = 34975
$ws.Range("N137").Value = -45175
$ws.Range("H139").Value = 45499.832
$ws.Range("J139").Value = 45499.832
$ws.Range("L139").Value = 45499.832
$ws.Range("N139").Value = -55779.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45992
$ws.Range("J16").Value = 45992
$ws.Range("L16").Value = 45992
$ws.Range("N16").Value = -46576
$ws.Range("H46").Value = 64953.09
$ws.Range("J46").Value = 64953.09
$ws.Range("L46").Value = 64953.09
$ws.Range("N46").Value = -65415.09
$ws.Range("H108").Value = 45626
$ws.Range("J108").Value = 45626
$ws.Range("L108").Value = 45626
$ws.Range("N108").Value = -53306
$ws.Range("H110").Value = 48636
$ws.Range("J110").Value = 48636
$ws.Range("L110").Value = 48636
$ws.Range("N110").Value = -56816
$ws.Range("H116").Value = 49680
$ws.Range("J116").Value = 49680
$ws.Range("L116").Value = 49680
$ws.Range("N116").Value = -58858
$ws.Range("H117").Value = 38421.4
$ws.Range("J117").Value = 38421.4
$ws.Range("L117").Value = 38421.4
$ws.Range("N117").Value = -47599.4
$ws.Range("H122").Value = 3572222.2
$ws.Range("I122").Value = 4082411.2
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 12247233.6
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -12244783.6
$ws.Range("N122").Value = -7600
$ws.Range("H125").Value = 35297.332
$ws.Range("J125").Value = 35297.332
$ws.Range("L125").Value = 35297.332
$ws.Range("N125").Value = -45137.332
$ws.Range("H131").Value = 50707
$ws.Range("J131").Value = 50707
$ws.Range("L131").Value = 50707
$ws.Range("N131").Value = -60787
$ws.Range("H132").Value = 1344.9286
$ws.Range("I132").Value = 1157.5333
$ws.Range("J132").Value = 2111.5454
$ws.Range("K132").Value = 3472.5999
$ws.Range("L132").Value = 6334.6362
$ws.Range("M132").Value = -942.5999000000002
$ws.Range("N132").Value = -11394.6362
$ws.Range("H134").Value = 64953.09
$ws.Range("J134").Value = 64953.09
$ws.Range("L134").Value = 194859.27
$ws.Range("N134").Value = -199929.27
$ws.Range("H138").Value = 46166.668
$ws.Range("J138").Value = 46166.668
$ws.Range("L138").Value = 46166.668
$ws.Range("N138").Value = -56446.668
$ws.Range("H139").Value = 38443
$ws.Range("I139").Value = 60958
$ws.Range("J139").Value = 33940
$ws.Range("K139").Value = 60958
$ws.Range("L139").Value = 33940
$ws.Range("N139").Value = -44220
